# Update database: shift yearly columns left (drop oldest 1396/12 period,
# add newest 1401/12 period) and change read_price algorithm values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Header row 8: financial period labels ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Header row 9: publish dates ---
$ws.Range("D9").Value = "1399-03-21 (8)"
$ws.Range("E9").Value = "1400-03-02 (8)"
$ws.Range("F9").Value = "1401-03-08 (8)"
$ws.Range("G9").Value = "1402-02-28 (7)"
# H9 is the bare date-looking string "1402-02-28"; Excel's text-to-date
# autoconvert would otherwise turn it into a date serial number, so force
# Text format just long enough to take the literal string, then restore the
# original General format/style (copied from a neighbouring untouched cell)
# so the cell's style stays identical to before.
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-28"
$ws.Range("G9").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4122) | Out-Null

# --- Row 11: فروش ---
$ws.Range("D11").Value = 25617
$ws.Range("E11").Value = 29989
$ws.Range("F11").Value = 26191
$ws.Range("G11").Value = 27409
$ws.Range("H11").Value = 35012

# --- Row 12: بهای تمام شده کالای فروش رفته ---
$ws.Range("D12").Value = -11759
$ws.Range("E12").Value = -14297
$ws.Range("F12").Value = -11451
$ws.Range("G12").Value = -15858
$ws.Range("H12").Value = -17813

# --- Row 13: سود (زیان) ناخالص ---
$ws.Range("D13").Value = 13857
$ws.Range("E13").Value = 15692
$ws.Range("F13").Value = 14740
$ws.Range("G13").Value = 11552
$ws.Range("H13").Value = 17199

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی ---
$ws.Range("D14").Value = -2667
$ws.Range("E14").Value = -1377
$ws.Range("F14").Value = -856
$ws.Range("G14").Value = -1658
$ws.Range("H14").Value = -1653

# --- Row 15: هزینه کاهش ارزش دریافتنی ها (هزینه استثنایی) ---
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = -202

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 1203
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 296
$ws.Range("G16").Value = -98
$ws.Range("H16").Value = -121

# --- Row 17: سود (زیان) عملیاتی ---
$ws.Range("D17").Value = 12393
$ws.Range("E17").Value = 14335
$ws.Range("F17").Value = 14180
$ws.Range("G17").Value = 9796
$ws.Range("H17").Value = 15224

# --- Row 18: هزینه های مالی ---
$ws.Range("D18").Value = -2390
$ws.Range("E18").Value = -2048
$ws.Range("F18").Value = -1369
$ws.Range("G18").Value = -2117
$ws.Range("H18").Value = -2867

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = -127
$ws.Range("E19").Value = 258
$ws.Range("F19").Value = 352
$ws.Range("G19").Value = 535
$ws.Range("H19").Value = 272

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 9876
$ws.Range("E20").Value = 12545
$ws.Range("F20").Value = 13163
$ws.Range("G20").Value = 8214
$ws.Range("H20").Value = 12629

# --- Row 21: مالیات ---
$ws.Range("D21").Value = -2319
$ws.Range("E21").Value = -2789
$ws.Range("F21").Value = -2154
$ws.Range("G21").Value = -1375
$ws.Range("H21").Value = -1861

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 7558
$ws.Range("E22").Value = 9756
$ws.Range("F22").Value = 11009
$ws.Range("G22").Value = 6839
$ws.Range("H22").Value = 10768

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (unchanged "-") ---
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"

# --- Row 24: سود (زیان) خالص ---
$ws.Range("D24").Value = 7558
$ws.Range("E24").Value = 9756
$ws.Range("F24").Value = 11009
$ws.Range("G24").Value = 6839
$ws.Range("H24").Value = 10768

# --- Row 25: سود هر سهم پس از کسر مالیات (unchanged zeros) ---
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

# --- Row 26: سرمایه ---
$ws.Range("D26").Value = 10182
$ws.Range("E26").Value = 8029
$ws.Range("F26").Value = 6767
$ws.Range("G26").Value = 9210
$ws.Range("H26").Value = 10286

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه (unchanged zeros) ---
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
